$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns for each coin row
# with the latest scraped values (rows 41/42 also swapped rank).

$ws.Range("D2").Value = "28.135.83"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.876.66"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'313.69"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.5128"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("D8").Value = "'0.3908"
$ws.Range("E8").Value = "  +2.48%  "
$ws.Range("D9").Value = "'0.08332"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").Value = "'1.120"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("D11").Value = "'41.44"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'6.219"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'20.63"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "1.873.25"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "'7.253"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "'0.00001100"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "'91.06"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "'0.06662"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "'17.76"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'6.030"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "28.176.00"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "'11.12"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "'2.254"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").Value = "2.088.98"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'2.489"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("D28").Value = "'159.40"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").Value = "'20.60"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").Value = "'125.09"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "'0.1062"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").Value = "'1.038"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").Value = "'5.841"
$ws.Range("E33").Value = "  +4.21%  "
$ws.Range("D34").Value = "'3.609"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").Value = "'9.609"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "'0.02450"
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").Value = "'0.06559"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "'0.2184"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").Value = "'1.198"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "'0.6489"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "'4.989"
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.226"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").Value = "'0.6136"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").Value = "'13.01"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "'3.670"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "'2.014"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").Value = "'1.231"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").Value = "'120.53"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'78.11"
$ws.Range("E51").Value = "  -1.31%  "
